# Update the "Förändrad" column (column C) for rows 2-19 from 45243 to 45244
# (i.e. advance the date value by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
